$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge "Hego" proofErr-wrapped run into the preceding text run in the
#    "parrainage" paragraph:
#    "...l'association étudiante " + [Hego] + " " -> "...l'association étudiante Hego "
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Replacement.ClearFormatting()
$rng1.Find.Execute("étudiante Hego ", $true, $false, $false, $false, $false, $true, 1, $false, "étudiante Hego ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Same merge for the "Le bureau des étudiants Hego Berria" paragraph.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Replacement.ClearFormatting()
$rng2.Find.Execute("étudiants Hego ", $true, $false, $false, $false, $false, $true, 1, $false, "étudiants Hego ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the "Saisie longue (identique à la saisie courte, mais plus
#    long)" list paragraph together with one of the two blank paragraphs
#    that follow it (before the "Brouillon algorithme" heading).
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Execute("Saisie longue (identique à la saisie courte, mais plus long)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$saisieLonguePara = $rng3.Paragraphs.Item(1)

$allParas = $d.Paragraphs
$startPos = $saisieLonguePara.Range.Start
$targetIndex = 0
$i = 1
foreach ($p in $allParas) {
  if ($p.Range.Start -eq $startPos) {
    $targetIndex = $i
  }
  $i = $i + 1
}
$pStart = $allParas.Item($targetIndex)
$pAfterBlank = $allParas.Item($targetIndex + 2)
$delRange = $d.Range($pStart.Range.Start, $pAfterBlank.Range.Start)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 4) Remove the hidden "_GoBack" bookmark around "On supprime les secondes
#    années qui ne veulent pas de filleuls".
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks.Item("_GoBack").Delete()
}

$d.Save()
